# Refresh crypto "Price" (D) and "Volume(1h)" (E) figures for the coinranking
# symbol list, per the scheduled GitHub Actions scrape (Tue Feb 7 10:31:51 UTC 2023).
# Values are stored as plain text (matching the sheet's existing inline-string
# convention), so each target cell is forced to Text format before the write and
# reset to the default "Normal" style afterwards (no visual/format side effects).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "330.82" },
    @{ Cell = "E2"; Value = "1.27%" },
    @{ Cell = "D3"; Value = "44.35" },
    @{ Cell = "E3"; Value = "-0.42%" },
    @{ Cell = "D4"; Value = "5.467" },
    @{ Cell = "E4"; Value = "-2.45%" },
    @{ Cell = "D5"; Value = "0.08033" },
    @{ Cell = "E5"; Value = "-0.44%" },
    @{ Cell = "D6"; Value = "1.987" },
    @{ Cell = "E6"; Value = "4.91%" },
    @{ Cell = "D7"; Value = "0.9539" },
    @{ Cell = "E7"; Value = "0.75%" },
    @{ Cell = "E8"; Value = "-3.95%" },
    @{ Cell = "D9"; Value = "0.1140" },
    @{ Cell = "E9"; Value = "-1.73%" },
    @{ Cell = "D10"; Value = "0.1896" },
    @{ Cell = "E10"; Value = "2.65%" },
    @{ Cell = "D11"; Value = "10.66" },
    @{ Cell = "E11"; Value = "26.37%" },
    @{ Cell = "D12"; Value = "0.09951" },
    @{ Cell = "E12"; Value = "0.75%" },
    @{ Cell = "D13"; Value = "0.04839" },
    @{ Cell = "E13"; Value = "14.38%" },
    @{ Cell = "D15"; Value = "0.001280" },
    @{ Cell = "E15"; Value = "-0.67%" },
    @{ Cell = "D16"; Value = "0.04081" },
    @{ Cell = "E16"; Value = "-3.26%" },
    @{ Cell = "D17"; Value = "0.005991" },
    @{ Cell = "E17"; Value = "0.57%" },
    @{ Cell = "D18"; Value = "3.367" },
    @{ Cell = "E18"; Value = "-6.73%" },
    @{ Cell = "D19"; Value = "4.392" },
    @{ Cell = "E19"; Value = "2.24%" },
    @{ Cell = "E20"; Value = "-2.07%" },
    @{ Cell = "D21"; Value = "0.1397" },
    @{ Cell = "E21"; Value = "1.91%" },
    @{ Cell = "D22"; Value = "0.2501" },
    @{ Cell = "E22"; Value = "-5.73%" },
    @{ Cell = "D23"; Value = "0.001273" },
    @{ Cell = "E23"; Value = "1.98%" },
    @{ Cell = "D24"; Value = "0.004361" },
    @{ Cell = "E24"; Value = "-3.26%" },
    @{ Cell = "E25"; Value = "-5.09%" },
    @{ Cell = "E26"; Value = "-6.35%" },
    @{ Cell = "D38"; Value = "0.02612" },
    @{ Cell = "E38"; Value = "-1.11%" },
    @{ Cell = "D39"; Value = "0.05809" },
    @{ Cell = "E39"; Value = "6.20%" },
    @{ Cell = "D40"; Value = "0.007553" },
    @{ Cell = "E40"; Value = "-0.93%" },
    @{ Cell = "D41"; Value = "0.1405" },
    @{ Cell = "E41"; Value = "0.60%" },
    @{ Cell = "D42"; Value = "0.007345" },
    @{ Cell = "E42"; Value = "0.00%" },
    @{ Cell = "E43"; Value = "-1.90%" },
    @{ Cell = "D44"; Value = "0.008859" },
    @{ Cell = "E44"; Value = "0.27%" },
    @{ Cell = "D45"; Value = "0.00006982" },
    @{ Cell = "E45"; Value = "0.64%" },
    @{ Cell = "E46"; Value = "-0.22%" },
    @{ Cell = "E48"; Value = "55.15%" },
    @{ Cell = "D49"; Value = "0.003551" },
    @{ Cell = "E49"; Value = "-4.19%" },
    @{ Cell = "E50"; Value = "-0.22%" },
    @{ Cell = "E51"; Value = "-0.22%" },
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"     # force Text so numeric-looking strings are not coerced
    $c.Value = $u.Value
    $c.Style = "Normal"       # drop back to the default style (no stray "@" format left on the cell)
}
